# Applies the "add correct results and results to debug of domapriori" edit:
#  1. Re-orders the object lists inside the rule descriptions on the "Reguły" sheet.
#  2. Re-labels / rotates the values on the "Walidacja krzyżowa" sheet so that the
#     accuracy / not_classified / correct / f1_score rows show the right figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Reguły" (Rules) sheet - reorder the object id lists quoted in each rule
# ---------------------------------------------------------------------------
$rules = $wb.Worksheets.Item("Reguły")

$rules.Range("B2").Value = "(attempts >=  3.0) => (class <= 1) ['a7', 'a23', 'a1', 'a13', 'a3']"
$rules.Range("B3").Value = "(pregnancy <=  0.0) & (age >=  32.0) & (frozen_embryos <=  4.0) & (cleavage_stage >=  5.0) => (class <= 1) ['a3', 'a15', 'a21', 'a22', 'a2']"
$rules.Range("B4").Value = "(sperm >=  3.0) => (class <= 1) ['a22', 'a25']"
$rules.Range("B5").Value = "(age >=  42.0) => (class <= 1) ['a14', 'a3']"
$rules.Range("B6").Value = "(age <=  31.0) & (attempts <=  1.0) => (class >= 2) ['a24', 'a11', 'a9', 'a25', 'a12']"
$rules.Range("B7").Value = "(frozen_embryos >=  8.0) & (sperm <=  1.0) => (class >= 2) ['a16', 'a6']"

# ---------------------------------------------------------------------------
# 2. "Walidacja krzyżowa" (Cross validation) sheet - rotate labels/values
# ---------------------------------------------------------------------------
$cv = $wb.Worksheets.Item("Walidacja krzyżowa")

$cv.Range("A1").Value = "accuracy"
$cv.Range("B1").Value = 0.4

$cv.Range("A2").Value = "not_classified"
$cv.Range("B2").Value = 0.44

$cv.Range("A3").Value = "correct"
$cv.Range("B3").Value = 0.7142857142857143

$cv.Range("A4").Value = "f1_score"
$cv.Range("B4").Value = 0.4952380952380952
